$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# S6: total post-processing time changed
$ws.Range("S6").Value = 1664

# B8..B26: the "name" column values were replaced with plain numeric-looking
# text labels (kept as text, not numbers)
$textCells = @{
    "B8"  = "146"
    "B9"  = "78"
    "B10" = "149"
    "B11" = "82"
    "B12" = "152"
    "B13" = "86"
    "B14" = "155"
    "B15" = "90"
    "B16" = "158"
    "B17" = "94"
    "B18" = "161"
    "B19" = "121"
    "B20" = "164"
    "B21" = "167"
    "B22" = "170"
    "B24" = "173"
    "B25" = "140"
    "B26" = "176"
}

foreach ($addr in $textCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$addr]
    $cell.Style = "Normal"
}

# S32: grand total time changed
$ws.Range("S32").Value = 28307

# Page margins reset to Excel's standard defaults
$ws.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.TopMargin = $excel.InchesToPoints(1)
$ws.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$ws.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$ws.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)
